$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 edits ---
# E6: "4" -> cleared (cell becomes empty)
$ws.Range("E6").ClearContents()

# G6: "6" -> "no pero si"
$ws.Range("G6").Value = "no pero si"

# --- Row 8 edits: clear C8, E8, F8, H8 entirely (D8 and G8 stay as-is) ---
$ws.Range("C8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("F8").ClearContents()
$ws.Range("H8").ClearContents()

# --- New row 9 ---
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "67555"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "3"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "a"

# --- New row 10 ---
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "12345"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "1"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "2"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1"
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = "5"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "COSINES"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "1"
